$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 0.57069473370060453
$ws.Range("C2").Value = 2.0273389758469595
$ws.Range("D2").Value = 0.45776957773598148
$ws.Range("E2").Value = 0.92486361928044003

$ws.Range("B3").Value = 0.91898135611724596
$ws.Range("C3").Value = 0.94677402519216591
$ws.Range("D3").Value = 0.67348917426278099
$ws.Range("E3").Value = 0.88001589369985722

$ws.Range("B1:E3").Select()
